$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "What is coronavirus / covid 19?"
$ws.Range("B2").Value = "Coronavirus is a type of virus which may cause illness in animals or humans, like the common cold, or more severe illnesses. covid 19 is a type of coronavirus which was recently discovered in China in 2019, and has spread to different parts of the world including India."
$ws.Range("B3").Value = "The most common symptoms of covid 19 are fever, tiredness, and dry cough. Some patients may have body pain, blocked nose, runny nose, sore throat or diarrhea. These symptoms are usually mild and begin gradually. Some people become infected but don’t develop any symptoms and don't feel unwell. Most people (about 80%) recover from the disease without needing special treatment. Around 1 out of every 6 people who gets covid 19 becomes seriously ill and develops difficulty breathing. Older people, and those with underlying medical problems like high blood pressure, heart problems or diabetes, are more likely to develop serious illness. People with fever, cough and difficulty breathing should seek medical attention."
$ws.Range("B4").Value = "People can catch covid 19 from others who have the virus. The disease can spread from person to person through small droplets from the nose or mouth which are spread when a person with covid 19 coughs or exhales. These droplets land on objects and surfaces around the person. Other people then catch covid 19 by touching these objects or surfaces, then touching their eyes, nose or mouth. People can also catch covid 19 if they breathe in droplets from a person with covid 19 who coughs out or exhales droplets. This is why it is important to stay more than 1 meter (3 feet) away from a person who is sick."
$ws.Range("B5").Value = "<ul><li>Regularly and properly washing your hands with soap and water, or cleaning with an alcohol based hand sanitizer kills viruses that may be on your hands.</li><li> Maintain at least 1 metre or 3 feet distance between yourself and anyone who is coughing or sneezing. When someone coughs or sneezes they spray small liquid droplets from their nose or mouth which may contain virus. If you are too close, you can breathe in the droplets, including the covid 19 virus if the person coughing has the disease.</li><li>Avoid touching eyes, nose and mouth.</li><li> Cover your mouth and nose with your bent elbow or tissue when you cough or sneeze. Then dispose of the used tissue immediately.</li><li>Stay home if you feel unwell. If you have a fever, cough and difficulty breathing, seek medical attention and call in advance.</li><li>Follow the directions of your local health authority.</li></ul>"

$ws.Range("G12").Select()